# Remove 'completed' items.
#
# Two items under the "MemoryMgr" heading have been completed and are
# removed in their entirety (including their paragraph mark, so the
# list collapses correctly around them):
#   - "Use Boost.FunctionTypes in remote function caller to detect calling
#      convention. (Also use TMP to detect number of args and their types etc)."
#   - "Improve genericity of parameter passing in remote function caller."

$d = $word.ActiveDocument

function Remove-ParagraphByText($text) {
    $rng = $d.Content
    $found = $rng.Find.Execute($text, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if ($found) {
        # Extend the found range by one character so the trailing paragraph
        # mark is included too, then delete via a fresh Range so the whole
        # paragraph (not just its text) disappears and the list collapses
        # around the gap, rather than leaving a blank line behind.
        $start = $rng.Start
        $delRange = $d.Range($start, $rng.End + 1)
        $delRange.Delete()
        return $start
    }
    return -1
}

$editPoint = Remove-ParagraphByText("Use Boost.FunctionTypes in remote function caller to detect calling convention. (Also use TMP to detect number of args and their types etc).")
$null = Remove-ParagraphByText("Improve genericity of parameter passing in remote function caller.")

# Word tracks the location of the most recent edit with the hidden
# '_GoBack' bookmark. Since this edit collapsed the deleted paragraph into
# the one that now follows it, move the bookmark there to reflect where the
# edit actually took place.
if ($editPoint -ge 0) {
    try {
        $d.Bookmarks.Item("_GoBack").Delete()
    } catch {
    }
    $editRange = $d.Range($editPoint, $editPoint)
    $d.Bookmarks.Add("_GoBack", $editRange)
}
